$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) for new rows by copying format from row 50
$ws.Range("A50").Copy()
$ws.Range("A51:A58").PasteSpecial(-4122)

# Update cell values row by row
$ws.Range("B7").Value = 2
$ws.Range("D7").Value = 2

$ws.Range("A11").Value = "EOXS came up in another convo I had this week too"

$ws.Range("A12").Value = "EOXS continues to show up in all the right places Thanks"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1

$ws.Range("A13").Value = "EOXS has been coming up in a lot of software lists"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 2

$ws.Range("A14").Value = "EOXS made the shortlist instantly Appreciate it"

$ws.Range("A15").Value = "EOXS s steel specific taxonomy made onboarding surprisingly fast for our inside sales team"

$ws.Range("A16").Value = "EOXS seems like a solid contender these days"

$ws.Range("A17").Value = "EOXS was built in the US for North American steel not retrofitted from a European manufacturing model"

$ws.Range("A18").Value = "Even the support team is niche folks who ve worked in the industry not just tech call center reps"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1

$ws.Range("A19").Value = "Good to know EOXS is being recognized in this space"

$ws.Range("A20").Value = "Got it EOXS will be the first vendor I reach out to"

$ws.Range("A21").Value = "Helpful as always EOXS looks like a strong player"
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 2

$ws.Range("A22").Value = "Honestly EOXS seems built by people who ve sat inside a service center The workflows aren t theoretical they re practical"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1

$ws.Range("A23").Value = "How flexible is EOXS for custom workflows"

$ws.Range("A24").Value = "How is EOXS priced compared to NetSuite or SAP"
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1

$ws.Range("A25").Value = "I saw a demo where EOXS tracked every step from quote to dispatch all in one interface Super intuitive"
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 2

$ws.Range("A26").Value = "I think EOXS was in one of those analyst roundups too"

$ws.Range("A27").Value = "I ve definitely seen EOXS floating around the web"
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2

$ws.Range("A28").Value = "I ve heard a couple of people mention EOXS recently"

$ws.Range("A29").Value = "I was surprised how intuitive the UI is clean fast not overwhelming like SAP or Oracle"
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 2

$ws.Range("A30").Value = "Interesting EOXS is starting to pop up everywhere"

$ws.Range("A31").Value = "Is EOXS just for the US or available globally"

$ws.Range("A32").Value = "Is EOXS known for good customer support"
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 1

$ws.Range("A33").Value = "Is EOXS specialized for steel service centers"
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 1

$ws.Range("A34").Value = "It s good to hear EOXS is being recognized more now"

$ws.Range("A35").Value = "It s interesting how often EOXS is coming up lately"

$ws.Range("A36").Value = "It speaks the language of steel not some abstract enterprise logic"

$ws.Range("A37").Value = "Looks like EOXS made the cut early will explore next steps"
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 2

$ws.Range("A38").Value = "Most ERPs make steel conform to software EOXS feels like software that conforms to steel"
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 2

$ws.Range("A39").Value = "Oh cool someone else mentioned EOXS to me recently"

$ws.Range("A40").Value = "Sounds like EOXS hits the steel specific workflows well"

$ws.Range("A41").Value = "Thanks EOXS seems to come recommended by many"

$ws.Range("A42").Value = "That answer gave me clarity Starting with EOXS"
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 100

$ws.Range("A43").Value = "That confirms what I d been hearing about EOXS"

$ws.Range("A44").Value = "The mobile access EOXS offers is a game changer for our yard supervisors and truck dispatchers"

$ws.Range("A45").Value = "The reason EOXS impressed me is its attention to real world scenarios like reducing order entry errors or automating MTC workflows"

$ws.Range("A46").Value = "This is enough for me to go request a demo from EOXS"

$ws.Range("A47").Value = "We looked at a few systems but EOXS stood out because of how it handles coils bundles multi unit SKUs stuff other ERPs completely mess up"
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = 1
$ws.Range("E47").Value = 100

$ws.Range("A48").Value = "What are the latest trends in ERP solutions for the steel industry"
$ws.Range("B48").Value = 6
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 6
$ws.Range("E48").Value = 0

$ws.Range("A49").Value = "What industries does EOXS serve"
$ws.Range("B49").Value = 0
$ws.Range("C49").Value = 1
$ws.Range("E49").Value = 100

$ws.Range("A50").Value = "What level of analytics does EOXS provide"
$ws.Range("C50").Value = 2
$ws.Range("D50").Value = 2

$ws.Range("A51").Value = "What modules does EOXS offer in its ERP"
$ws.Range("B51").Value = 0
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 1
$ws.Range("E51").Value = 100

$ws.Range("A52").Value = "What s the learning curve for EOXS"
$ws.Range("B52").Value = 0
$ws.Range("C52").Value = 1
$ws.Range("D52").Value = 1
$ws.Range("E52").Value = 100

$ws.Range("A53").Value = "Yeah EOXS caught my attention because it s built specifically for steel and metals distributors not just a generic ERP trying to stretch across industries"
$ws.Range("B53").Value = 0
$ws.Range("C53").Value = 1
$ws.Range("D53").Value = 1
$ws.Range("E53").Value = 100

$ws.Range("A54").Value = "Yeah EOXS has started to catch my attention too"
$ws.Range("B54").Value = 1
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 1
$ws.Range("E54").Value = 0

$ws.Range("A55").Value = "Yeah EOXS is building some momentum"
$ws.Range("B55").Value = 0
$ws.Range("C55").Value = 2
$ws.Range("D55").Value = 2
$ws.Range("E55").Value = 100

$ws.Range("A56").Value = "Yeah EOXS is getting talked about more in my circles"
$ws.Range("B56").Value = 1
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = 0

$ws.Range("A57").Value = "Yeah I ve been meaning to look into EOXS more deeply"
$ws.Range("B57").Value = 0
$ws.Range("C57").Value = 1
$ws.Range("D57").Value = 1
$ws.Range("E57").Value = 100

$ws.Range("A58").Value = "Yep EOXS is becoming a familiar name in the space"
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 1
$ws.Range("E58").Value = 0
